$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- "Hárok1" view: scroll back to the top and select column D ---
$ws1.Range("D1:D1048576").Select() | Out-Null

# --- Add the new "Hárok2" worksheet right after "Hárok1" ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Hárok2"

# Column headers
$ws2.Range("A1").Value = "no-analyzers"
$ws2.Range("B1").Value = "bh-analyzers"

# Bulk-copy the already-computed "minutes" columns from "Hárok1" across:
#   column B (no-analyzer minutes, incl. the Trimmean on row 102) -> A
#   column D (bh-analyzer minutes,  incl. the Trimmean on row 102) -> B
$ws2.Range("A2:A102").Value = $ws1.Range("B2:B102").Value()
$ws2.Range("B2:B102").Value = $ws1.Range("D2:D102").Value()

# Overhead summary (ratio / seconds saved) - column B only, row 103 stays blank
$ws2.Range("B104").Value = $ws1.Range("D104").Value()
$ws2.Range("B105").Value = $ws1.Range("D105").Value()

# Column widths to roughly match the authored sheet (engine quantises to 1/6
# character-width steps, so these are the closest achievable values to the
# authored 17.42578125 / 13.140625)
$ws2.Columns.Item(1).ColumnWidth = 16.59
$ws2.Columns.Item(2).ColumnWidth = 12.3

# Leave the new sheet selected on F8, matching the saved view
$ws2.Range("F8").Select() | Out-Null
